# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (column G) values computed/regenerated for each data row (rows 2-20)
$kValues = @{
    2  = 4
    3  = 7
    4  = 6
    5  = 5
    6  = 8
    7  = 6
    8  = 8
    9  = 8
    10 = 10
    11 = 10
    12 = 5
    13 = 3
    14 = 7
    15 = 4
    16 = 7
    17 = 11
    18 = 5
    19 = 4
    20 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
